$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 6: table's style id changes from the custom "Table_0" style to the
#    built-in table style {AADD72BC-3906-4830-9A2E-7EE958A8B4BF}.
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$tblShape = $s6.Shapes.Item(2)
$tbl = $tblShape.Table
$tbl.ApplyStyle("{AADD72BC-3906-4830-9A2E-7EE958A8B4BF}", $false)

# ---------------------------------------------------------------------------
# 2) Theme swap: the deck's active theme (shared by the slide master / the
#    presentation) switches from the custom "Integral" palette to the
#    stock "Office Theme" palette.
# ---------------------------------------------------------------------------
$master = $p.SlideMaster
$colors = $master.ColorScheme

# index -> (scheme slot, new colour)
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
# 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$colors.Colors(1).RGB  = 0x000000   # dk1      -> 000000
$colors.Colors(2).RGB  = 0xFFFFFF   # lt1      -> FFFFFF
$colors.Colors(3).RGB  = 0x6A5444   # dk2      -> 44546A
$colors.Colors(4).RGB  = 0xE6E6E7   # lt2      -> E7E6E6
$colors.Colors(5).RGB  = 0xD59B5B   # accent1  -> 5B9BD5
$colors.Colors(6).RGB  = 0x317DED   # accent2  -> ED7D31
$colors.Colors(7).RGB  = 0xA5A5A5   # accent3  -> A5A5A5
$colors.Colors(8).RGB  = 0x00C0FF   # accent4  -> FFC000
$colors.Colors(9).RGB  = 0xC47244   # accent5  -> 4472C4
$colors.Colors(10).RGB = 0x47AD70   # accent6  -> 70AD47
$colors.Colors(11).RGB = 0xC16305   # hlink    -> 0563C1
$colors.Colors(12).RGB = 0x724F95   # folHlink -> 954F72
